# Auto-generated Excel COM-interop edit script
# Applies a scheduled market-data refresh to the H:N "profit" columns
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (per commit message:
# "chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 600
$ws.Range("I21").Value = 600
$ws.Range("K21").Value = 600
$ws.Range("M21").Value = -132

$ws.Range("H23").Value = 600
$ws.Range("I23").Value = 600
$ws.Range("K23").Value = 600
$ws.Range("M23").Value = -366

$ws.Range("H112").Value = 5566808.5
$ws.Range("J112").Value = 6343482
$ws.Range("L112").Value = 19030446
$ws.Range("N112").Value = -19032662

$ws.Range("H133").Value = 11333.333
$ws.Range("J133").Value = 11333.333
$ws.Range("L133").Value = 11333.333
$ws.Range("N133").Value = -21453.333

$ws.Range("H138").Value = 12562464
$ws.Range("I138").Value = 1501399.2
$ws.Range("J138").Value = 333333340
$ws.Range("K138").Value = 4504197.6
$ws.Range("L138").Value = 1000000020
$ws.Range("M138").Value = -4499057.6
$ws.Range("N138").Value = -1000010300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 712
$ws.Range("I45").Value = 712
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 712
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -335
$ws.Range("N45").ClearContents()

$ws.Range("H122").Value = 8229.666999999999
$ws.Range("I122").Value = 8229.666999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24689.001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -22239.001
$ws.Range("N122").ClearContents()

$ws.Range("H133").Value = 46379.8
$ws.Range("J133").Value = 46379.8
$ws.Range("L133").Value = 46379.8
$ws.Range("N133").Value = -51439.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3496.5
$ws.Range("J20").Value = 2490
$ws.Range("L20").Value = 2490
$ws.Range("N20").Value = -2984

$ws.Range("H94").Value = 1261.2106
$ws.Range("I94").Value = 1031.4375
$ws.Range("J94").Value = 2486.6667
$ws.Range("K94").Value = 1031.4375
$ws.Range("L94").Value = 2486.6667
$ws.Range("M94").Value = -580.4375
$ws.Range("N94").Value = -3388.6667

$ws.Range("H107").Value = 846.25
$ws.Range("I107").Value = 801.41174
$ws.Range("K107").Value = 801.41174
$ws.Range("M107").Value = 1118.58826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2098.4243
$ws.Range("I31").Value = 1211.76
$ws.Range("J31").Value = 4869.25
$ws.Range("K31").Value = 1211.76
$ws.Range("L31").Value = 4869.25
$ws.Range("M31").Value = -916.76
$ws.Range("N31").Value = -5459.25

$ws.Range("H34").Value = 2098.4243
$ws.Range("I34").Value = 1211.76
$ws.Range("J34").Value = 4869.25
$ws.Range("K34").Value = 1211.76
$ws.Range("L34").Value = 4869.25
$ws.Range("M34").Value = -1009.76
$ws.Range("N34").Value = -5273.25

$ws.Range("H58").Value = 2452.8462
$ws.Range("I58").Value = 1450.9412
$ws.Range("J58").Value = 4345.3335
$ws.Range("K58").Value = 1450.9412
$ws.Range("L58").Value = 4345.3335
$ws.Range("M58").Value = -1247.9412
$ws.Range("N58").Value = -4751.3335

$ws.Range("J86").Value = 1691.5625
$ws.Range("L86").Value = 1691.5625
$ws.Range("N86").Value = -3937.5625

$ws.Range("J89").Value = 1691.5625
$ws.Range("L89").Value = 8457.8125
$ws.Range("N89").Value = -19689.8125

$ws.Range("H132").Value = 2354.5
$ws.Range("I132").Value = 1843.3572
$ws.Range("J132").Value = 5932.5
$ws.Range("K132").Value = 5530.071599999999
$ws.Range("L132").Value = 17797.5
$ws.Range("M132").Value = -3000.071599999999
$ws.Range("N132").Value = -22857.5

$ws.Range("H134").Value = 2522.9285
$ws.Range("I134").Value = 1690.7142
$ws.Range("J134").Value = 6684
$ws.Range("K134").Value = 5072.142599999999
$ws.Range("L134").Value = 20052
$ws.Range("M134").Value = -2537.142599999999
$ws.Range("N134").Value = -25122

$ws.Range("H136").Value = 2452.8462
$ws.Range("I136").Value = 1450.9412
$ws.Range("J136").Value = 4345.3335
$ws.Range("K136").Value = 4352.8236
$ws.Range("L136").Value = 13036.0005
$ws.Range("M136").Value = -1802.8236
$ws.Range("N136").Value = -18136.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1292.8889
$ws.Range("I97").Value = 1119.5714
$ws.Range("K97").Value = 1119.5714
$ws.Range("M97").Value = -623.5714

$ws.Range("H102").Value = 3195.7334
$ws.Range("I102").Value = 3317
$ws.Range("K102").Value = 3317
$ws.Range("M102").Value = -1695

$ws.Range("H132").Value = 3730.7058
$ws.Range("I132").Value = 3481.2593
$ws.Range("K132").Value = 10443.7779
$ws.Range("M132").Value = -7913.777900000001

$ws.Range("H138").Value = 64466.668
$ws.Range("J138").Value = 64466.668
$ws.Range("L138").Value = 64466.668
$ws.Range("N138").Value = -74746.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4349160
$ws.Range("I16").Value = 4763320
$ws.Range("K16").Value = 4763320
$ws.Range("M16").Value = -4763150

$ws.Range("H22").Value = 12555.889
$ws.Range("I22").Value = 1699.5
$ws.Range("J22").Value = 15657.714
$ws.Range("K22").Value = 1699.5
$ws.Range("L22").Value = 15657.714
$ws.Range("M22").Value = -1404.5
$ws.Range("N22").Value = -16247.714

$ws.Range("H27").Value = 12555.889
$ws.Range("I27").Value = 1699.5
$ws.Range("J27").Value = 15657.714
$ws.Range("K27").Value = 1699.5
$ws.Range("L27").Value = 15657.714
$ws.Range("M27").Value = -1592.5
$ws.Range("N27").Value = -15871.714

$ws.Range("H40").Value = 3198.4666
$ws.Range("I40").Value = 1639.9166
$ws.Range("K40").Value = 1639.9166
$ws.Range("M40").Value = -1503.9166

$ws.Range("H122").Value = 3357.762
$ws.Range("I122").Value = 2027
$ws.Range("J122").Value = 3670.8823
$ws.Range("K122").Value = 6081
$ws.Range("L122").Value = 11012.6469
$ws.Range("M122").Value = -3631
$ws.Range("N122").Value = -15912.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 466809.66
$ws.Range("J46").Value = 466809.66
$ws.Range("L46").Value = 466809.66
$ws.Range("N46").Value = -467271.66

$ws.Range("H126").Value = 46359.184
$ws.Range("I126").Value = 67441.13
$ws.Range("J126").Value = 1183.5714
$ws.Range("K126").Value = 202323.39
$ws.Range("L126").Value = 3550.7142
$ws.Range("M126").Value = -199853.39
$ws.Range("N126").Value = -8490.7142

$ws.Range("H132").Value = 9805956
$ws.Range("I132").Value = 12501922
$ws.Range("J132").Value = 2445.4546
$ws.Range("K132").Value = 37505766
$ws.Range("L132").Value = 7336.3638
$ws.Range("M132").Value = -37503236
$ws.Range("N132").Value = -12396.3638

$ws.Range("H134").Value = 466809.66
$ws.Range("J134").Value = 466809.66
$ws.Range("L134").Value = 1400428.98
$ws.Range("N134").Value = -1405498.98
